$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E13 status from "CV Sent" to "1st Interview"
$ws.Range("E13").Value = "1st Interview"

# Add new row 15 with data
$ws.Range("A15").Value = 867
$ws.Range("B15").Value = "CyCognito"
$ws.Range("C15").Value = "CSM UK"
$ws.Range("D15").Value = "Gary M."
$ws.Range("E15").Value = "CV Sent"
